$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.16%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.105"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.31%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05665"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.81%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.469"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.45%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8222"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.19%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8437"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.46%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1324"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.65%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06926"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.57%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02881"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.24%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09386"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.23%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001514"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.14%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04119"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-12.05%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-93.92%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006204"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.76%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.513"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.98%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.86%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.311"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "9.11%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03150"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.35%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1291"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.16%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.555"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.27%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.07%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.32%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004451"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.73%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009800"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2.11%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "3.47%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03674"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.06%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006050"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.28%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1052"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.61%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002300"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.41%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009302"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.26%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005315"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.02%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-15.83%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002567"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "24.30%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.02%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"
